# Adds a new paragraph to the "Content Placeholder 6" shape on the
# "Navigating the MFDs" slide, describing the MFDs' Option Select Buttons
# (OSBs). That shape's text currently ends with:
#   "The AMPCD can display the same page as an MFD."
# and we append a further paragraph after it.

$p = $ppt.ActivePresentation

$marker = "The AMPCD can display the same page as an MFD."

$targetSlide = $null
$targetShape = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $shapeText = $shape.TextFrame.TextRange.Text
            if ($shapeText -like "*$marker*") {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

# Fall back to the known location if the search above failed to find it.
if ($targetShape -eq $null) {
    $targetSlide = $p.Slides.Item(9)
    $targetShape = $targetSlide.Shapes.Item(2)
}

$textRange = $targetShape.TextFrame.TextRange

$firstRunText = "The MDFs each have 20 Option Select Buttons (OSBs), numbered from 1 to 20. OSB 1 is the bottom-most button on the left side of the MFD, with the numbers increasing clockwise around the MFD. Accordingly, OSB 20 is the left-most OSB on the "
$secondRunText = "bottom row of buttons."

# Insert a new paragraph after the existing text, then append the second
# part of the sentence separately so that it ends up as its own run within
# the new paragraph (matching the two differently-touched runs in the
# target content).
$inserted = $textRange.InsertAfter("`r" + $firstRunText)
[void]$inserted.InsertAfter($secondRunText)
